$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.870.20"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "3.084.82"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "237.69"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "608.58"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").Value = "1.11"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "0.383"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "0.798"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("D11").Value = "3.083.68"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "0.196"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "94.400.93"
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("D14").Value = "0.0000241"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "33.65"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").Value = "5.34"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").Value = "3.661.08"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "3.080.27"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "3.59"
$ws.Range("E19").Value = "  -4.86%  "
$ws.Range("D20").Value = "14.35"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").Value = "5.72"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").Value = "445.48"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "8.89"
$ws.Range("E23").Value = "  -4.06%  "
$ws.Range("D24").Value = "0.0000193"
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "5.53"
$ws.Range("E25").Value = "  -3.14%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "84.78"
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "11.73"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "3.245.32"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.135"
$ws.Range("E30").Value = "  +2.50%  "
$ws.Range("D31").Value = "0.243"
$ws.Range("E31").Value = "  +3.06%  "
$ws.Range("D32").Value = "0.174"
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "0.997"
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "8.96"
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").Value = "7.43"
$ws.Range("E35").Value = "  -5.15%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "25.59"
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("E37").Value = "  -4.18%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "490.04"
$ws.Range("E38").Value = "  +4.09%  "
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("B40").Value = "MantraDAO"
$ws.Range("C40").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D40").Value = "3.84"
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "24.07"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "0.433"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "1.24"
$ws.Range("E43").Value = "  -3.23%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "3.16"
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "161.44"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "0.679"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.83"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "0.000279"
$ws.Range("E49").Value = "  +13.64%  "
$ws.Range("D50").Value = "43.71"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("E51").Value = "  +0.07%  "
